# kalyankar to pusapati done
# Fill in grading marks ("Total Points") and a grading comment for the
# toString() method rows of the Customer Class and Product Class sections,
# then leave the selection where the grader was last working (F14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Customer Class section (rows 3-6) ---
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "(-1) for wrong output in toString method"

# --- Product Class section (rows 10-14) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "(-1) for wrong output in toString method"

# Move the active selection to where grading left off.
$ws.Range("F14").Select()
